# TC04_CDS_phs002293_SampleTumorStatus_Normal.xlsx
# "Fixed CDS phs002293 TC04 and TC05"
#
# The FilesTab query (row 4 / "FilesTab", column "TabQuery" -> cell B4)
# had its trailing LIMIT clause commented out ("--LIMIT 100;"). Uncomment
# it so the query actually limits the result set again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Cells.Item(4, 2)
$query = $cell.Value2
$fixedQuery = $query.Replace("--LIMIT 100;", "LIMIT 100;")
$cell.Value2 = $fixedQuery

# Match the saved cursor position left in the workbook (selection on B3).
$ws.Range("B3").Select()

$wb.Save()
